# JackieK-WorkLog.xlsx - "Add files via upload"
# Adds a new work-log entry (row 45) for 2025-04-01, describing another
# validation test run with TryHackMe, and normalizes the formatting of the
# previous "latest entry" row (44) back to the regular (non-highlighted)
# style now that row 45 is the newest entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Carry row 44's current ("most recent entry") formatting down onto the
#    new row 45 before we touch row 44 itself.
$ws.Range("A44:C44").Copy()
$ws.Range("A45:C45").PasteSpecial(-4122)

# 2) Reset row 44 (B:C) back to the plain/regular formatting used by the
#    other completed rows (it is no longer the newest entry).
$ws.Range("B41:C41").Copy()
$ws.Range("B44:C44").PasteSpecial(-4122)

# 3) Fill in the new log entry's data.
$ws.Range("A45").Value = 45748
$ws.Range("B45").Value = 4
$ws.Range("C45").Value = "Ran another validation test with ""TryHackMe"" and updated final report"

# 4) Match the saved selection/cursor position recorded in the workbook.
$ws.Range("C48").Select() | Out-Null
